$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5060.5
$ws.Range("I6").Value = 5060.5
$ws.Range("K6").Value = 15181.5
$ws.Range("M6").Value = -15069.5
$ws.Range("H9").Value = 12155.2
$ws.Range("I9").Value = 12155.2
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 12155.2
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -11986.2
$ws.Range("N9").Value = $null
$ws.Range("H12").Value = 350
$ws.Range("I12").Value = 350
$ws.Range("K12").Value = 350
$ws.Range("M12").Value = -180
$ws.Range("H21").Value = 144.66667
$ws.Range("I21").Value = 144.66667
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 144.66667
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 323.33333
$ws.Range("N21").Value = $null
$ws.Range("H23").Value = 144.66667
$ws.Range("I23").Value = 144.66667
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 144.66667
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 89.33332999999999
$ws.Range("N23").Value = $null
$ws.Range("H29").Value = 1980
$ws.Range("J29").Value = 2412.5
$ws.Range("L29").Value = 7237.5
$ws.Range("N29").Value = -7799.5
$ws.Range("H38").Value = 381.69232
$ws.Range("I38").Value = 106.111115
$ws.Range("J38").Value = 1001.75
$ws.Range("K38").Value = 318.333345
$ws.Range("L38").Value = 3005.25
$ws.Range("M38").Value = 53.66665499999999
$ws.Range("N38").Value = -3749.25
$ws.Range("H58").Value = 1332.9166
$ws.Range("I58").Value = 263
$ws.Range("J58").Value = 2402.8333
$ws.Range("K58").Value = 789
$ws.Range("L58").Value = 7208.499899999999
$ws.Range("M58").Value = -639
$ws.Range("N58").Value = -7508.499899999999
$ws.Range("H86").Value = 4112.4375
$ws.Range("I86").Value = 3957
$ws.Range("J86").Value = 4233.3335
$ws.Range("K86").Value = 3957
$ws.Range("L86").Value = 4233.3335
$ws.Range("M86").Value = -2834
$ws.Range("N86").Value = -6479.3335
$ws.Range("H87").Value = 53000
$ws.Range("J87").Value = 53000
$ws.Range("L87").Value = 53000
$ws.Range("N87").Value = -55496
$ws.Range("H89").Value = 4112.4375
$ws.Range("I89").Value = 3957
$ws.Range("J89").Value = 4233.3335
$ws.Range("K89").Value = 19785
$ws.Range("L89").Value = 21166.6675
$ws.Range("M89").Value = -14169
$ws.Range("N89").Value = -32398.6675
$ws.Range("H90").Value = 53000
$ws.Range("J90").Value = 53000
$ws.Range("L90").Value = 159000
$ws.Range("N90").Value = -171480
$ws.Range("H98").Value = 2123.3333
$ws.Range("I98").Value = 1731.4286
$ws.Range("K98").Value = 1731.4286
$ws.Range("M98").Value = -233.4286
$ws.Range("H122").Value = 2123.3333
$ws.Range("I122").Value = 1731.4286
$ws.Range("K122").Value = 5194.2858
$ws.Range("M122").Value = -2744.2858
$ws.Range("H132").Value = 7940418.5
$ws.Range("I132").Value = 8549240
$ws.Range("J132").Value = 25733.334
$ws.Range("K132").Value = 25647720
$ws.Range("L132").Value = 77200.00199999999
$ws.Range("M132").Value = -25645190
$ws.Range("N132").Value = -82260.00199999999
$ws.Range("H141").Value = 835.5
$ws.Range("I141").Value = 835.5
$ws.Range("K141").Value = 2506.5
$ws.Range("M141").Value = 2673.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 22680
$ws.Range("I31").Value = 5450
$ws.Range("K31").Value = 5450
$ws.Range("M31").Value = -5156
$ws.Range("H61").Value = 1440.5
$ws.Range("I61").Value = 1089
$ws.Range("K61").Value = 1089
$ws.Range("M61").Value = -877
$ws.Range("H74").Value = 1027.0454
$ws.Range("I74").Value = 820.0526
$ws.Range("J74").Value = 2338
$ws.Range("K74").Value = 820.0526
$ws.Range("L74").Value = 2338
$ws.Range("M74").Value = 53.94740000000002
$ws.Range("N74").Value = -4086
$ws.Range("H77").Value = 1027.0454
$ws.Range("I77").Value = 820.0526
$ws.Range("J77").Value = 2338
$ws.Range("K77").Value = 4100.263
$ws.Range("L77").Value = 11690
$ws.Range("M77").Value = 267.7370000000001
$ws.Range("N77").Value = -20426
$ws.Range("H136").Value = 1440.5
$ws.Range("I136").Value = 1089
$ws.Range("K136").Value = 3267
$ws.Range("M136").Value = -717

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 17416.5
$ws.Range("I102").Value = 8125
$ws.Range("J102").Value = 35999.5
$ws.Range("K102").Value = 8125
$ws.Range("L102").Value = 35999.5
$ws.Range("M102").Value = -4880
$ws.Range("N102").Value = -42489.5
$ws.Range("H134").Value = 7516.2856
$ws.Range("I134").Value = 1552.6875
$ws.Range("J134").Value = 26599.8
$ws.Range("K134").Value = 4658.0625
$ws.Range("L134").Value = 79799.39999999999
$ws.Range("M134").Value = -2123.0625
$ws.Range("N134").Value = -84869.39999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
$ws.Range("H99").Value = 1828.2941
$ws.Range("J99").Value = 1860.4286
$ws.Range("L99").Value = 1860.4286
$ws.Range("N99").Value = -4856.4286
$ws.Range("H126").Value = 1828.2941
$ws.Range("J126").Value = 1860.4286
$ws.Range("L126").Value = 5581.2858
$ws.Range("N126").Value = -10521.2858
$ws.Range("H132").Value = 3579.7273
$ws.Range("I132").Value = 3063.3333
$ws.Range("J132").Value = 4199.4
$ws.Range("K132").Value = 9189.999899999999
$ws.Range("L132").Value = 12598.2
$ws.Range("M132").Value = -6659.999899999999
$ws.Range("N132").Value = -17658.2
$ws.Range("H134").Value = 14493701
$ws.Range("I134").Value = 17544634
$ws.Range("J134").Value = 1775
$ws.Range("K134").Value = 52633902
$ws.Range("L134").Value = 5325
$ws.Range("M134").Value = -52631367
$ws.Range("N134").Value = -10395

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 849
$ws.Range("I122").Value = 849
$ws.Range("K122").Value = 7641
$ws.Range("M122").Value = -5191
$ws.Range("H132").Value = 1121.6
$ws.Range("I132").Value = 1102
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 9918
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -7388
$ws.Range("N132").Value = -15860

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5101.8335
$ws.Range("I80").Value = 2752.5
$ws.Range("J80").Value = 6276.5
$ws.Range("K80").Value = 2752.5
$ws.Range("L80").Value = 6276.5
$ws.Range("M80").Value = -1754.5
$ws.Range("N80").Value = -8272.5
$ws.Range("H83").Value = 5101.8335
$ws.Range("I83").Value = 2752.5
$ws.Range("J83").Value = 6276.5
$ws.Range("K83").Value = 13762.5
$ws.Range("L83").Value = 31382.5
$ws.Range("M83").Value = -8770.5
$ws.Range("N83").Value = -41366.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H113").Value = 3726.2856
$ws.Range("I113").Value = 1711
$ws.Range("J113").Value = 4845.8887
$ws.Range("K113").Value = 1711
$ws.Range("L113").Value = 4845.8887
$ws.Range("M113").Value = 459
$ws.Range("N113").Value = -9185.8887
$ws.Range("H122").Value = 251234.5
$ws.Range("I122").Value = 1481.4
$ws.Range("K122").Value = 4444.200000000001
$ws.Range("M122").Value = -1994.200000000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 14166.667
$ws.Range("J94").Value = 14166.667
$ws.Range("L94").Value = 14166.667
$ws.Range("N94").Value = -15518.667
$ws.Range("H136").Value = 1977.4546
$ws.Range("I136").Value = 1779
$ws.Range("K136").Value = 5337
$ws.Range("M136").Value = -2787

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 12030.6
$ws.Range("J51").Value = 12538.5
$ws.Range("L51").Value = 12538.5
$ws.Range("N51").Value = -13558.5
$ws.Range("H52").Value = 16588
$ws.Range("J52").Value = 16588
$ws.Range("L52").Value = 16588
$ws.Range("N52").Value = -17040
$ws.Range("H132").Value = 3572.9546
$ws.Range("I132").Value = 3562.9375
$ws.Range("J132").Value = 3599.6667
$ws.Range("K132").Value = 10688.8125
$ws.Range("L132").Value = 10799.0001
$ws.Range("M132").Value = -8158.8125
$ws.Range("N132").Value = -15859.0001
